$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally has two stacked header rows (row1 / row2) and a stray
# partial data row (row5: F5:H5 only) that doesn't belong to any indexed
# canton entry. Remove the stray row first (higher row index first so row
# numbers of rows above it are unaffected), then drop the two old header
# rows, and finally insert a single new combined header row at the top.

$ws.Rows.Item(5).Delete()   # stray row (F=3, G=5.3, H=5.3) - discarded
$ws.Rows.Item(2).Delete()   # old 2nd header row (units row)
$ws.Rows.Item(1).Delete()   # old 1st header row

$ws.Rows.Item(1).Insert()   # fresh blank row for the new combined header

# New left-hand identification columns
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

# Measurement columns keep the smaller 9pt Arial font used elsewhere in the
# table (same font as the data rows' header cells previously used).
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

$hdr = $ws.Range("F1:K1")
$hdr.Font.Name = "Arial"
$hdr.Font.Size = 9

$ws.Range("A4:K4").Select()
